$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "W:\Projects\הסעת המונים\מטרו\01_שלב ה\קבצי עבודה\תחזיות_דמוגרפיות\תחזיות_2050\התפלגות גילים\backround_files"
$ws.Range("C8").Value = "create_age_distribution_230719.ipynb"
